$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: log the extra session that day and append the note about getting sick
$ws.Range("B17").Value = 120
$ws.Range("C17").Value = "Further worked on request handling"
$ws.Range("D17").Value = "16:00-18:00; i got sick :("

# Row 18: new entry - reuse the "01.01.2024" text already used in column A,
# and copy D17's time-note style (text, numFmt) for the new note cell.
$ws.Range("A15").Copy()
$ws.Range("A18").PasteSpecial(-4163)

$ws.Range("D17").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "11:00-zeit"

$ws.Application.CutCopyMode = $false

# Move the selection to the last-edited cell, like the author did
$ws.Range("D18").Select()
